$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Partner: [Имя партнёра] – Insertion Sort" -> split "Partner: " / name /
#    " – Insertion Sort" across three separate runs.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Partner: [Имя партнёра] – Insertion Sort", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rng.Find.Found) {
    throw "Could not locate the Partner paragraph"
}
$p = $rng.Paragraphs(1)
$target = $d.Range($p.Range.Start, $p.Range.End - 1)
$target.InsertXML('<w:p><w:r><w:t xml:space="preserve">Partner: </w:t></w:r><w:r><w:t>Temirlan Almukhamedov</w:t></w:r><w:r><w:t xml:space="preserve"> – Insertion Sort</w:t></w:r></w:p>')

# ---------------------------------------------------------------------------
# 2) Merge the standalone "<w:r><w:br/></w:r>" run that precedes
#    "- The algorithm performs about n² / 2 comparisons..." into the run
#    carrying that sentence's text.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Explanation:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rng.Find.Found) {
    throw "Could not locate the Explanation paragraph"
}
$p = $rng.Paragraphs(1)
$target = $d.Range($p.Range.Start, $p.Range.End - 1)
$target.InsertXML('<w:p><w:r><w:t>Explanation:</w:t></w:r><w:r><w:br/><w:t>- The algorithm performs about n² / 2 comparisons in the average and worst cases.</w:t></w:r><w:r><w:br/><w:t>- The number of swaps is linear (O(n)) because each iteration moves at most two elements.</w:t></w:r><w:r><w:br/><w:t>- The algorithm is in-place and uses only a few auxiliary variables.</w:t></w:r></w:p>')

# ---------------------------------------------------------------------------
# 3) Merge "...educational purp" / "oses." back into a single run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("The measured results confirm the theoretical analysis:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rng.Find.Found) {
    throw "Could not locate the measured-results paragraph"
}
$p = $rng.Paragraphs(1)
$target = $d.Range($p.Range.Start, $p.Range.End - 1)
$target.InsertXML('<w:p><w:r><w:t>The measured results confirm the theoretical analysis:</w:t></w:r><w:r><w:br/><w:t>- Best Case (O(n)) — if the array is already sorted, the algorithm stops early.</w:t></w:r><w:r><w:br/><w:t>- Average/Worst Case (O(n²)) — for unsorted data, comparisons grow quadratically.</w:t></w:r><w:r><w:br/><w:t>- Optimization Effect: Early termination reduced time for nearly sorted arrays.</w:t></w:r><w:r><w:br/><w:t>- Space Efficiency: Algorithm is in-place with O(1) additional memory.</w:t></w:r><w:r><w:br/></w:r><w:r><w:br/><w:t>Despite the optimizations, Selection Sort remains inefficient for large datasets because of its O(n²) complexity. However, it is simple, predictable, and useful for small arrays or educational purposes.</w:t></w:r></w:p>')
